$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.902.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.910.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9979"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.57%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.66"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9982"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4995"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.76%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3806"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07287"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.99%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.31"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9088"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07638"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.880.84"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.474"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.54"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9981"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008726"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.75%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9981"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "27.943.43"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.46%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.131.96"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.15%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.602"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.60%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.89"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.862"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.219"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.09%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.64%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.98"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.904"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08988"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.198"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.810"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.231"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7800"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.635"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02084"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.060"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.090"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5541"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.05278"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.807"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "113.69"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.515"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1515"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.61"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4825"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9984"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.639"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.27"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06041"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.65%  "
